$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -18.74353603170066
$ws.Range("C2").Value = -18.74353603170066
$ws.Range("D2").Value = -18.74353603170066
$ws.Range("E2").Value = -18.74353603170066
$ws.Range("F2").Value = -18.74353603170066
$ws.Range("G2").Value = -18.74353603170066
$ws.Range("H2").Value = -18.74353603170066
$ws.Range("I2").Value = -18.74353603170066
$ws.Range("J2").Value = -18.74353603170066
$ws.Range("K2").Value = -18.74353603170066
$ws.Range("B3").Value = -18.74353603170066
$ws.Range("C3").Value = -18.74353603170066
$ws.Range("D3").Value = -18.74353603170066
$ws.Range("E3").Value = -18.74353603170066
$ws.Range("F3").Value = -18.74353603170066
$ws.Range("G3").Value = -18.74353603170066
$ws.Range("H3").Value = -18.74353603170066
$ws.Range("I3").Value = 2.513672812035318
$ws.Range("J3").Value = -18.74353603170066
$ws.Range("K3").Value = -18.74353603170066
$ws.Range("B4").Value = -18.74353603170066
$ws.Range("C4").Value = -18.74353603170066
$ws.Range("D4").Value = 3.146652900091476
$ws.Range("E4").Value = -18.74353603170066
$ws.Range("F4").Value = 2.496801708922677
$ws.Range("G4").Value = -18.74353603170066
$ws.Range("H4").Value = 1.905787693366359
$ws.Range("I4").Value = -18.74353603170066
$ws.Range("J4").Value = 2.391071401200536
$ws.Range("K4").Value = -18.74353603170066
$ws.Range("B5").Value = -18.74353603170066
$ws.Range("C5").Value = -18.74353603170066
$ws.Range("D5").Value = -18.74353603170066
$ws.Range("E5").Value = -18.74353603170066
$ws.Range("F5").Value = -18.74353603170066
$ws.Range("G5").Value = 2.167956169078088
$ws.Range("H5").Value = -18.74353603170066
$ws.Range("I5").Value = -18.74353603170066
$ws.Range("J5").Value = -18.74353603170066
$ws.Range("K5").Value = -18.74353603170066
$ws.Range("B6").Value = -18.74353603170066
$ws.Range("C6").Value = -18.74353603170066
$ws.Range("D6").Value = -18.74353603170066
$ws.Range("E6").Value = -18.74353603170066
$ws.Range("F6").Value = -18.74353603170066
$ws.Range("G6").Value = -18.74353603170066
$ws.Range("H6").Value = -18.74353603170066
$ws.Range("I6").Value = -18.74353603170066
$ws.Range("J6").Value = -18.74353603170066
$ws.Range("K6").Value = -18.74353603170066
$ws.Range("B7").Value = 2.996541566166514
$ws.Range("C7").Value = -18.74353603170066
$ws.Range("D7").Value = -18.74353603170066
$ws.Range("E7").Value = -18.74353603170066
$ws.Range("F7").Value = -18.74353603170066
$ws.Range("G7").Value = -18.74353603170066
$ws.Range("H7").Value = -18.74353603170066
$ws.Range("I7").Value = -18.74353603170066
$ws.Range("J7").Value = -18.74353603170066
$ws.Range("K7").Value = -18.74353603170066
$ws.Range("B8").Value = -18.74353603170066
$ws.Range("C8").Value = -18.74353603170066
$ws.Range("D8").Value = -18.74353603170066
$ws.Range("E8").Value = 3.000046429011439
$ws.Range("F8").Value = -18.74353603170066
$ws.Range("G8").Value = -18.74353603170066
$ws.Range("H8").Value = -18.74353603170066
$ws.Range("I8").Value = -18.74353603170066
$ws.Range("J8").Value = -18.74353603170066
$ws.Range("K8").Value = -18.74353603170066
$ws.Range("B9").Value = 3.587258603230743
$ws.Range("C9").Value = -18.74353603170066
$ws.Range("D9").Value = -18.74353603170066
$ws.Range("E9").Value = -18.74353603170066
$ws.Range("F9").Value = -18.74353603170066
$ws.Range("G9").Value = -18.74353603170066
$ws.Range("H9").Value = -18.74353603170066
$ws.Range("I9").Value = -18.74353603170066
$ws.Range("J9").Value = -18.74353603170066
$ws.Range("K9").Value = -18.74353603170066
$ws.Range("B10").Value = -18.74353603170066
$ws.Range("C10").Value = -18.74353603170066
$ws.Range("D10").Value = -18.74353603170066
$ws.Range("E10").Value = -18.74353603170066
$ws.Range("F10").Value = -18.74353603170066
$ws.Range("G10").Value = -18.74353603170066
$ws.Range("H10").Value = -18.74353603170066
$ws.Range("I10").Value = 1.534360884754535
$ws.Range("J10").Value = -18.74353603170066
$ws.Range("K10").Value = 2.245587477218284
$ws.Range("B11").Value = -18.74353603170066
$ws.Range("C11").Value = -18.74353603170066
$ws.Range("D11").Value = -18.74353603170066
$ws.Range("E11").Value = 2.205898161605235
$ws.Range("F11").Value = -18.74353603170066
$ws.Range("G11").Value = 2.607571968462443
$ws.Range("H11").Value = -18.74353603170066
$ws.Range("I11").Value = -18.74353603170066
$ws.Range("J11").Value = -18.74353603170066
$ws.Range("K11").Value = 1.437272761742195
$ws.Range("B12").Value = -18.74353603170066
$ws.Range("C12").Value = -18.74353603170066
$ws.Range("D12").Value = -18.74353603170066
$ws.Range("E12").Value = -18.74353603170066
$ws.Range("F12").Value = -18.74353603170066
$ws.Range("G12").Value = -18.74353603170066
$ws.Range("H12").Value = -18.74353603170066
$ws.Range("I12").Value = -18.74353603170066
$ws.Range("J12").Value = -18.74353603170066
$ws.Range("K12").Value = -18.74353603170066
$ws.Range("B13").Value = -18.74353603170066
$ws.Range("C13").Value = -18.74353603170066
$ws.Range("D13").Value = -18.74353603170066
$ws.Range("E13").Value = 1.684956025949134
$ws.Range("F13").Value = -18.74353603170066
$ws.Range("G13").Value = -18.74353603170066
$ws.Range("H13").Value = -18.74353603170066
$ws.Range("I13").Value = -18.74353603170066
$ws.Range("J13").Value = 2.221586631233621
$ws.Range("K13").Value = 1.629240837786083
$ws.Range("B14").Value = -18.74353603170066
$ws.Range("C14").Value = -18.74353603170066
$ws.Range("D14").Value = 1.601974254653563
$ws.Range("E14").Value = -18.74353603170066
$ws.Range("F14").Value = -18.74353603170066
$ws.Range("G14").Value = -18.74353603170066
$ws.Range("H14").Value = -18.74353603170066
$ws.Range("I14").Value = -18.74353603170066
$ws.Range("J14").Value = -18.74353603170066
$ws.Range("K14").Value = 2.088616139106907
$ws.Range("B15").Value = -18.74353603170066
$ws.Range("C15").Value = -18.74353603170066
$ws.Range("D15").Value = 0.5719518879344757
$ws.Range("E15").Value = -18.74353603170066
$ws.Range("F15").Value = -18.74353603170066
$ws.Range("G15").Value = -18.74353603170066
$ws.Range("H15").Value = -18.74353603170066
$ws.Range("I15").Value = -18.74353603170066
$ws.Range("J15").Value = -18.74353603170066
$ws.Range("K15").Value = -18.74353603170066
$ws.Range("B16").Value = -18.74353603170066
$ws.Range("C16").Value = -18.74353603170066
$ws.Range("D16").Value = -18.74353603170066
$ws.Range("E16").Value = -18.74353603170066
$ws.Range("F16").Value = -18.74353603170066
$ws.Range("G16").Value = -18.74353603170066
$ws.Range("H16").Value = -18.74353603170066
$ws.Range("I16").Value = -18.74353603170066
$ws.Range("J16").Value = 2.305475585281771
$ws.Range("K16").Value = -18.74353603170066
$ws.Range("B17").Value = -18.74353603170066
$ws.Range("C17").Value = -18.74353603170066
$ws.Range("D17").Value = 0.7144239158687942
$ws.Range("E17").Value = -18.74353603170066
$ws.Range("F17").Value = -18.74353603170066
$ws.Range("G17").Value = -18.74353603170066
$ws.Range("H17").Value = 0.7519930812664806
$ws.Range("I17").Value = 0.9456719991497485
$ws.Range("J17").Value = 1.28348576904085
$ws.Range("K17").Value = -18.74353603170066
$ws.Range("B18").Value = -18.74353603170066
$ws.Range("C18").Value = -18.74353603170066
$ws.Range("D18").Value = -18.74353603170066
$ws.Range("E18").Value = -18.74353603170066
$ws.Range("F18").Value = -18.74353603170066
$ws.Range("G18").Value = -18.74353603170066
$ws.Range("H18").Value = 0.6996038035775048
$ws.Range("I18").Value = 1.112635663082103
$ws.Range("J18").Value = 1.439819484526486
$ws.Range("K18").Value = -18.74353603170066
$ws.Range("B19").Value = -18.74353603170066
$ws.Range("C19").Value = -18.74353603170066
$ws.Range("D19").Value = 1.287477222337968
$ws.Range("E19").Value = -18.74353603170066
$ws.Range("F19").Value = -18.74353603170066
$ws.Range("G19").Value = -18.74353603170066
$ws.Range("H19").Value = 1.782123364952106
$ws.Range("I19").Value = 1.971112805883118
$ws.Range("J19").Value = -18.74353603170066
$ws.Range("K19").Value = -18.74353603170066
$ws.Range("B20").Value = -18.74353603170066
$ws.Range("C20").Value = 4.321924972174559
$ws.Range("D20").Value = 1.3448446216315
$ws.Range("E20").Value = -18.74353603170066
$ws.Range("F20").Value = 3.84354488575394
$ws.Range("G20").Value = -18.74353603170066
$ws.Range("H20").Value = 2.132412004120148
$ws.Range("I20").Value = 1.758496440161036
$ws.Range("J20").Value = -18.74353603170066
$ws.Range("K20").Value = 2.379254687275369
$ws.Range("B21").Value = -18.74353603170066
$ws.Range("C21").Value = -18.74353603170066
$ws.Range("D21").Value = -18.74353603170066
$ws.Range("E21").Value = 2.060318387163913
$ws.Range("F21").Value = -18.74353603170066
$ws.Range("G21").Value = 3.234398216547647
$ws.Range("H21").Value = 2.356379971475427
$ws.Range("I21").Value = -18.74353603170066
$ws.Range("J21").Value = -18.74353603170066
$ws.Range("K21").Value = -18.74353603170066
